# Weekly refresh: insert 3 new data rows (one new date, 3 quality grades)
# at the top of the "Acelga" price block, pushing the existing rows down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows at row 233, shifting existing rows 233-324 down to 236-327.
$ws.Rows("233:235").Insert()

# --- New row 233 : Extra quality, new date 44468 ---
$ws.Cells.Item(233, 1).Value = 9
$ws.Cells.Item(233, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(233, 3).Value = "Metropolitana"
$ws.Cells.Item(233, 4).Value = 44468
$ws.Cells.Item(233, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(233, 5).Value = 13
$ws.Cells.Item(233, 6).Value = 100112009
$ws.Cells.Item(233, 7).Value = "Acelga"
$ws.Cells.Item(233, 8).Value = "Sin especificar"
$ws.Cells.Item(233, 9).Value = "Extra"
$ws.Cells.Item(233, 10).Value = 18
$ws.Cells.Item(233, 11).Value = 12000
$ws.Cells.Item(233, 12).Value = 12000
$ws.Cells.Item(233, 13).Value = 12000
$ws.Cells.Item(233, 14).Value = "`$/docena de atados"
$ws.Cells.Item(233, 15).Value = "Región Metropolitana"
$ws.Cells.Item(233, 16).Value = 4000
$ws.Cells.Item(233, 17).Value = 3
$ws.Cells.Item(233, 18).Value = "Hortaliza"

# --- New row 234 : Primera quality, new date 44468 ---
$ws.Cells.Item(234, 1).Value = 9
$ws.Cells.Item(234, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(234, 3).Value = "Metropolitana"
$ws.Cells.Item(234, 4).Value = 44468
$ws.Cells.Item(234, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(234, 5).Value = 13
$ws.Cells.Item(234, 6).Value = 100112009
$ws.Cells.Item(234, 7).Value = "Acelga"
$ws.Cells.Item(234, 8).Value = "Sin especificar"
$ws.Cells.Item(234, 9).Value = "Primera"
$ws.Cells.Item(234, 10).Value = 52
$ws.Cells.Item(234, 11).Value = 10000
$ws.Cells.Item(234, 12).Value = 11000
$ws.Cells.Item(234, 13).Value = 10500
$ws.Cells.Item(234, 14).Value = "`$/docena de atados"
$ws.Cells.Item(234, 15).Value = "Región Metropolitana"
$ws.Cells.Item(234, 16).Value = 3500
$ws.Cells.Item(234, 17).Value = 3
$ws.Cells.Item(234, 18).Value = "Hortaliza"

# --- New row 235 : Segunda quality, new date 44468 ---
$ws.Cells.Item(235, 1).Value = 9
$ws.Cells.Item(235, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(235, 3).Value = "Metropolitana"
$ws.Cells.Item(235, 4).Value = 44468
$ws.Cells.Item(235, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(235, 5).Value = 13
$ws.Cells.Item(235, 6).Value = 100112009
$ws.Cells.Item(235, 7).Value = "Acelga"
$ws.Cells.Item(235, 8).Value = "Sin especificar"
$ws.Cells.Item(235, 9).Value = "Segunda"
$ws.Cells.Item(235, 10).Value = 34
$ws.Cells.Item(235, 11).Value = 8000
$ws.Cells.Item(235, 12).Value = 9000
$ws.Cells.Item(235, 13).Value = 8500
$ws.Cells.Item(235, 14).Value = "`$/docena de atados"
$ws.Cells.Item(235, 15).Value = "Región Metropolitana"
$ws.Cells.Item(235, 16).Value = 2833
$ws.Cells.Item(235, 17).Value = 3
$ws.Cells.Item(235, 18).Value = "Hortaliza"
